$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E) entirely, shifting the
# remaining columns (F:K) one position to the left.
$ws.Columns.Item(5).Delete()
